$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reactions")

# --- New column H: "Objective proportion" -----------------------------
# Clone the header formatting used by the other E1:G1 header cells (bold
# font on a light fill) onto H1, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Objective proportion"

# New data column values for rows 2-6.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 1

# --- E2:E6 become real boolean cells -----------------------------------
$ws.Range("E2").Value = $false
$ws.Range("E3").Value = $false
$ws.Range("E4").Value = $true
$ws.Range("E5").Value = $false

# E6 previously had a different (non-boolean) style than E2:E5; align its
# formatting with the rest of the column before writing the boolean value.
$ws.Range("E6").Style = $ws.Range("E2").Style
$ws.Range("E6").Value = $false

# --- Selection / active sheet ------------------------------------------
# The edited workbook now has the "Reactions" sheet active (instead of
# "Rate laws"), with H1:H6 selected.
$ws.Activate() | Out-Null
$ws.Range("H1:H6").Select() | Out-Null
